$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet currently ends with a "data" block (rows 2-103) followed by a single
# footer/note row at row 104 (empty date cell + a shared-string footnote in column B).
# New data for 2020-05-08 needs to be inserted as the new last data row (104),
# and the footer/note row needs to shift down to row 105.

# 1) Duplicate the existing footer row (104) down into a new row 105, pushing
#    nothing else around (there is no data below row 104).
$ws.Rows.Item(104).Copy() | Out-Null
$ws.Rows.Item(105).Insert(-4121) | Out-Null

# 2) Row 104 is now free of content (still holds the old footer's formatting).
#    Bring in the number formatting used by the rest of the data rows (copy it
#    from the previous data row, 103) before writing the new values.
$ws.Range("A103:E103").Copy() | Out-Null
$ws.Range("A104:E104").PasteSpecial(-4122) | Out-Null

# 3) Write the new day's figures (2020-05-08) into row 104.
$ws.Range("A104").Value = 43959
$ws.Range("B104").Value = 304
$ws.Range("C104").Value = 35007
$ws.Range("D104").Value = 99
$ws.Range("E104").Value = 7134

# 4) Move the active selection to the new last data row, matching where the
#    sheet's cursor ends up after appending a row.
$ws.Range("A104").Select() | Out-Null

# 5) The printed area now needs to include the extra row.
$wb.Names.Item(1).RefersTo = "=" + $ws.Name + "!`$A`$1:`$E`$106"

Write-Host "Appended 2020-05-08 row; footer moved to row 105; Print_Area extended."
